$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.243.02'
$ws.Range('E2').Value = '  -3.12%  '

$ws.Range('D3').Value = '3.753.31'
$ws.Range('E3').Value = '  -0.46%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.67'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -3.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.88'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -6.03%  '

$ws.Range('D7').Value = '3.750.84'
$ws.Range('E7').Value = '  -0.42%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.158'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -5.34%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.18'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -4.77%  '

$ws.Range('E12').Value = '  -4.86%  '

$ws.Range('E13').Value = '  -4.56%  '

$ws.Range('E14').Value = '  -5.25%  '

$ws.Range('D15').Value = '4.392.31'
$ws.Range('E15').Value = '  -0.17%  '

$ws.Range('D16').Value = '3.756.72'
$ws.Range('E16').Value = '  -0.20%  '

$ws.Range('D17').Value = '67.323.18'
$ws.Range('E17').Value = '  -3.08%  '

$ws.Range('E18').Value = '  -3.85%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.13'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -5.56%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.33'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +5.54%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '487.96'
$ws.Range('D21').NumberFormat = 'General'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.23'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -1.88%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.729'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.27%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.83'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -1.56%  '

$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000144'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.35'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -5.61%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.14'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -5.73%  '

$ws.Range('E28').Value = '  -4.83%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.94'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -2.37%  '

$ws.Range('E31').Value = '  -5.84%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.17'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +3.75%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.74'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -4.54%  '

$ws.Range('E34').Value = '  -6.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('E36').Value = '  -3.68%  '

$ws.Range('E37').Value = '  -5.37%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.132'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -5.45%  '

$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '458.80'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -1.74%  '

$ws.Range('E40').Value = '  -4.38%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '49.19'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -1.37%  '

$ws.Range('E42').Value = '  -4.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.81'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -7.10%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.31'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -3.20%  '

$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.13'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -9.40%  '

$ws.Range('D47').Value = '2.809.26'
$ws.Range('E47').Value = '  -4.79%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '140.82'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +1.11%  '

$ws.Range('E49').Value = '  -4.30%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.06'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +9.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.50'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -6.86%  '
